$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range('D2').NumberFormat = '@'
$ws.Range('D2').Value = '68.296.96'
$ws.Range('D2').Style = 'Normal'
$ws.Range('E2').Value = '  +1.78%  '
$ws.Range('D3').NumberFormat = '@'
$ws.Range('D3').Value = '3.896.42'
$ws.Range('D3').Style = 'Normal'
$ws.Range('E3').Value = '  +1.04%  '
$ws.Range('E4').Value = '  +0.19%  '
$ws.Range('D5').NumberFormat = '@'
$ws.Range('D5').Value = '479.92'
$ws.Range('D5').Style = 'Normal'
$ws.Range('E5').Value = '  +2.44%  '
$ws.Range('D6').NumberFormat = '@'
$ws.Range('D6').Value = '144.60'
$ws.Range('D6').Style = 'Normal'
$ws.Range('E6').Value = '  -0.13%  '
$ws.Range('E7').Value = '  -1.88%  '
$ws.Range('E8').Value = '  -0.10%  '
$ws.Range('D9').NumberFormat = '@'
$ws.Range('D9').Value = '0.722'
$ws.Range('D9').Style = 'Normal'
$ws.Range('E9').Value = '  -2.92%  '
$ws.Range('E10').Value = '  +7.88%  '
$ws.Range('D11').NumberFormat = '@'
$ws.Range('D11').Value = '0.0000353'
$ws.Range('D11').Style = 'Normal'
$ws.Range('E11').Value = '  +14.00%  '
$ws.Range('D12').NumberFormat = '@'
$ws.Range('D12').Value = '42.60'
$ws.Range('D12').Style = 'Normal'
$ws.Range('E12').Value = '  -1.67%  '
$ws.Range('D13').NumberFormat = '@'
$ws.Range('D13').Value = '10.60'
$ws.Range('D13').Style = 'Normal'
$ws.Range('E13').Value = '  +2.12%  '
$ws.Range('D14').NumberFormat = '@'
$ws.Range('D14').Value = '4.525.97'
$ws.Range('D14').Style = 'Normal'
$ws.Range('E14').Value = '  +1.36%  '
$ws.Range('D15').NumberFormat = '@'
$ws.Range('D15').Value = '14.61'
$ws.Range('D15').Style = 'Normal'
$ws.Range('E15').Value = '  -1.35%  '
$ws.Range('D16').NumberFormat = '@'
$ws.Range('D16').Value = '3.905.16'
$ws.Range('D16').Style = 'Normal'
$ws.Range('E16').Value = '  +1.02%  '
$ws.Range('E17').Value = '  -0.32%  '
$ws.Range('D18').NumberFormat = '@'
$ws.Range('D18').Value = '19.70'
$ws.Range('D18').Style = 'Normal'
$ws.Range('E18').Value = '  -1.74%  '
$ws.Range('E19').Value = '  -3.26%  '
$ws.Range('D20').NumberFormat = '@'
$ws.Range('D20').Value = '68.348.80'
$ws.Range('D20').Style = 'Normal'
$ws.Range('E20').Value = '  +1.60%  '
$ws.Range('D21').NumberFormat = '@'
$ws.Range('D21').Value = '435.98'
$ws.Range('D21').Style = 'Normal'
$ws.Range('E21').Value = '  +0.50%  '
$ws.Range('D22').NumberFormat = '@'
$ws.Range('D22').Value = '14.74'
$ws.Range('D22').Style = 'Normal'
$ws.Range('E22').Value = '  -1.27%  '
$ws.Range('E23').Value = '  +0.72%  '
$ws.Range('D24').NumberFormat = '@'
$ws.Range('D24').Value = '88.11'
$ws.Range('D24').Style = 'Normal'
$ws.Range('E24').Value = '  -0.77%  '
$ws.Range('D25').NumberFormat = '@'
$ws.Range('D25').Value = '11.74'
$ws.Range('D25').Style = 'Normal'
$ws.Range('E25').Value = '  +18.37%  '
$ws.Range('E26').Value = '  -0.82%  '
$ws.Range('D27').NumberFormat = '@'
$ws.Range('D27').Value = '10.42'
$ws.Range('D27').Style = 'Normal'
$ws.Range('E27').Value = '  +3.09%  '
$ws.Range('D28').NumberFormat = '@'
$ws.Range('D28').Value = '38.03'
$ws.Range('D28').Style = 'Normal'
$ws.Range('E28').Value = '  +0.27%  '
$ws.Range('D29').NumberFormat = '@'
$ws.Range('D29').Value = '5.81'
$ws.Range('D29').Style = 'Normal'
$ws.Range('E29').Value = '  +4.57%  '
$ws.Range('D30').NumberFormat = '@'
$ws.Range('D30').Value = '709.96'
$ws.Range('D30').Style = 'Normal'
$ws.Range('E30').Value = '  -2.47%  '
$ws.Range('D31').NumberFormat = '@'
$ws.Range('D31').Value = '0.130'
$ws.Range('D31').Style = 'Normal'
$ws.Range('E31').Value = '  -2.53%  '
$ws.Range('D32').NumberFormat = '@'
$ws.Range('D32').Value = '13.34'
$ws.Range('D32').Style = 'Normal'
$ws.Range('E32').Value = '  -3.48%  '
$ws.Range('D33').NumberFormat = '@'
$ws.Range('D33').Value = '2.86'
$ws.Range('D33').Style = 'Normal'
$ws.Range('E33').Value = '  +2.42%  '
$ws.Range('D34').NumberFormat = '@'
$ws.Range('D34').Value = '0.0₃0924'
$ws.Range('D34').Style = 'Normal'
$ws.Range('E34').Value = '  +36.76%  '
$ws.Range('D35').NumberFormat = '@'
$ws.Range('D35').Value = '41.55'
$ws.Range('D35').Style = 'Normal'
$ws.Range('E35').Value = '  -5.36%  '
$ws.Range('D36').NumberFormat = '@'
$ws.Range('D36').Value = '59.36'
$ws.Range('D36').Style = 'Normal'
$ws.Range('E36').Value = '  +2.47%  '
$ws.Range('D37').NumberFormat = '@'
$ws.Range('D37').Value = '5.68'
$ws.Range('D37').Style = 'Normal'
$ws.Range('E37').Value = '  +3.22%  '
$ws.Range('E38').Value = '  -6.13%  '
$ws.Range('D40').NumberFormat = '@'
$ws.Range('D40').Value = '0.0474'
$ws.Range('D40').Style = 'Normal'
$ws.Range('E40').Value = '  -2.00%  '
$ws.Range('D41').NumberFormat = '@'
$ws.Range('D41').Value = '3.07'
$ws.Range('D41').Style = 'Normal'
$ws.Range('E41').Value = '  +10.95%  '
$ws.Range('E42').Value = '  +8.02%  '
$ws.Range('E43').Value = '  +2.65%  '
$ws.Range('B44').Value = 'TheGraph'
$ws.Range('C44').Value = 'https://coinranking.com/coin/qhd1biQ7M+thegraph-grt'
$ws.Range('D44').NumberFormat = '@'
$ws.Range('D44').Value = '0.340'
$ws.Range('D44').Style = 'Normal'
$ws.Range('E44').Value = '  -2.05%  '
$ws.Range('B45').Value = 'Stellar'
$ws.Range('C45').Value = 'https://coinranking.com/coin/f3iaFeCKEmkaZ+stellar-xlm'
$ws.Range('D45').NumberFormat = '@'
$ws.Range('D45').Value = '0.141'
$ws.Range('D45').Style = 'Normal'
$ws.Range('E45').Value = '  -0.19%  '
$ws.Range('E46').Value = '  +0.11%  '
$ws.Range('E47').Value = '  -1.16%  '
$ws.Range('E48').Value = '  -0.57%  '
$ws.Range('D49').NumberFormat = '@'
$ws.Range('D49').Value = '145.84'
$ws.Range('D49').Style = 'Normal'
$ws.Range('E49').Value = '  +1.19%  '
$ws.Range('D50').NumberFormat = '@'
$ws.Range('D50').Value = '3.12'
$ws.Range('D50').Style = 'Normal'
$ws.Range('E50').Value = '  -4.78%  '
$ws.Range('E51').Value = '  -1.99%  '
